# Print request logsheet: add Equipment Tray Holder 1 + 2 entries for
# 26-07-2018, and fill in newly-completed "Date Completed" values for
# three previously-open rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "July 2018" sheet, tab index 2 (0-based)

# --- Update "Date Completed" for rows that have now finished printing ---
# (match the centered style already used throughout column B)
$ws.Range("B6").Value  = "24-07-2018"   # Spine Box Case
$ws.Range("B6").HorizontalAlignment = -4108

$ws.Range("B9").Value  = "26-07-2018"   # FOCUS lax 100 Scale
$ws.Range("B9").HorizontalAlignment = -4108

$ws.Range("B16").Value = "25-07-2018"   # Equipment Tray Clip
$ws.Range("B16").HorizontalAlignment = -4108

# --- Append the two new "Equipment Tray Holder" print requests ---
$ws.Range("A17").Value = "26-07-2018"
$ws.Range("C17").Value = "Equipment Tray Holder 1"
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = "Polylite"
$ws.Range("F17").Value = 4
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = 0.2
$ws.Range("I17").Value = "NA"

$ws.Range("A18").Value = "26-07-2018"
$ws.Range("C18").Value = "Equipment Tray Holder 2"
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = "Polylite"
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = 0.2
$ws.Range("I18").Value = "NA"

# --- Update the sheet's last selection to match the new extent ---
$ws.Range("B18").Select()
